$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update data rows 2-5 with new values ---
$row2 = @(45050.50694444445, 22.58, 15.542, 4.221, 47.493, 39.284, 17.769, 58.8, 27.341, 11.61, 17.881, 18.828, 19.728, 5.673, 17.67, 24.849, 14.79, 3.779, 2.46, 261.617, 49.202, 16.31, 32.642, 17.025, 2.109, 29.266, 14.407, 12.944, 15.145, 19.485, 3.64, 51.902, 9.071, 20.391)
$row3 = @(45050.51388888889, 11.53, 7.975, 1.653, 24.437, 20.228, 9.074, 37.778, 13.961, 5.962, 9.044, 9.839, 10.209, 2.901, 9.023, 12.69, 7.794, 1.582, 0.916, 130.032, 25.348, 8.329000000000001, 16.716, 8.962999999999999, 1.025, 17.826, 7.357, 6.717, 7.851, 10.212, 1.294, 34.195, 4.578, 10.413)
$row4 = @(45050.52083333334, 12.491, 8.930999999999999, 1.192, 26.753, 22.136, 9.83, 38.568, 15.125, 6.58, 9.897, 10.818, 11.258, 3.14, 9.775, 13.807, 8.375, 1.05, 0.695, 141.456, 27.338, 9.023, 18.179, 9.766999999999999, 1.154, 18.432, 7.97, 7.183, 8.417999999999999, 11.28, 0.784, 34.723, 5.01, 11.28)
$row5 = @(45050.52777777778, 14.89, 10.85, 1.06, 32.08, 26.54, 11.72, 45.27, 18.03, 7.94, 11.89, 12.97, 13.55, 3.74, 11.65, 16.53, 9.890000000000001, 0.82, 0.66, 170.07, 32.58, 10.76, 21.81, 11.65, 1.42, 21.89, 9.5, 8.48, 9.960000000000001, 13.56, 0.5600000000000001, 40.87, 6.02, 13.45)

$rows = @($row2, $row3, $row4, $row5)

for ($ri = 0; $ri -lt $rows.Length; $ri++) {
    $rowData = $rows[$ri]
    $excelRow = $ri + 2
    for ($ci = 0; $ci -lt $rowData.Length; $ci++) {
        $ws.Cells.Item($excelRow, $ci + 1).Value = $rowData[$ci]
    }
}

# --- Remove the now-obsolete row 6 (original 5th data row), shrinking the range to A1:AH5 ---
$ws.Rows.Item(6).Delete()

# --- Widen specific columns from 7 to 8 characters ---
$colsToWiden = @(2, 3, 7, 11, 12, 13, 16, 24, 27, 28, 29, 30, 34)
foreach ($colIndex in $colsToWiden) {
    $ws.Columns.Item($colIndex).ColumnWidth = 7.16
}
